$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..AI to D..AJ),
# then insert a new (mostly blank) row before row 4 (shifts old rows 4-7 to 5-8).
$ws.Columns("C").Insert()
$ws.Rows("4").Insert()

# Give the new "Project Number 2" column roughly the same width as the
# "Project Number" column (B) instead of the default width.
$ws.Range("C1").ColumnWidth = $ws.Range("B1").ColumnWidth

# New header in C1 and duplicate of the Project Number value in C2.
$ws.Range("C1").Value = "Project Number 2"
$ws.Range("C2").Value = $ws.Range("B2").Value2

# Match the selection left behind in the source workbook.
$ws.Range("C2").Select() | Out-Null
